# For subsector/variable rows 3, 4, 5, 6 and 9 on the active sheet, the
# per-year trajectory (columns K through AS, i.e. years 1-35) is
# overwritten so every year takes on the same value as year 0
# (column J) for that row - effectively flattening the trajectory to a
# constant, consistent with the commit updating the rail freight/
# passenger diesel and electric efficiency trajectories.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(3, 4, 5, 6, 9)

foreach ($row in $rows) {
    $baseValue = $ws.Range("J$row").Value2
    $ws.Range("K$row`:AS$row").Value2 = $baseValue
}
